$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 82 (pushing the existing register_14..register_18 rows
# down to 84..88, where they become register_15..register_19).
$ws.Range("A82:A83").EntireRow.Insert()

# Copy the formatting (style indices) of a row that already has the exact
# "blank bit-field row" look (cols B..L styled 27,8,8,9,13,13,14,12,46,46,47)
# onto the two freshly-inserted rows.
$ws.Range("B74:L74").Copy()
$ws.Range("B82:L82").PasteSpecial(-4122)
$ws.Range("B83:L83").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows(82).RowHeight = 15.75
$ws.Rows(83).RowHeight = 15.75

# Row 82: new register_14 register row (offset 0x64, counter type, size 8)
$ws.Range("B82").Value = "register_14"
$ws.Range("C82").Value = "0x64"
$ws.Range("D82").Value = ""
$ws.Range("E82").Value = ""
$ws.Range("F82").Value = "name"
$ws.Range("G82").Value = "8"
$ws.Range("H82").Value = "counter"
$ws.Range("I82").Value = 0

# Row 83: the counter bit field belonging to register_14
$ws.Range("B83").Value = ""
$ws.Range("C83").Value = ""
$ws.Range("D83").Value = ""
$ws.Range("E83").Value = ""
$ws.Range("F83").Value = "bit_field_1"
$ws.Range("G83").Value = "8"
$ws.Range("H83").Value = "counter"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = "register_3.bit_field_3"

# Rows 84..88 already hold the old register_14..register_18 rows (shifted
# down by the insert above) with all formatting/values intact; only their
# register names need to be renumbered by one.
$ws.Range("B84").Value = "register_15"
$ws.Range("B85").Value = "register_16"
$ws.Range("B86").Value = "register_17"
$ws.Range("B87").Value = "register_18"
$ws.Range("B88").Value = "register_19"

# Update the view state to match the saved workbook.
$ws.Range("C83").Select()
